$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export")

# Row 5 (004472404 DILSON 54650.18) is replaced in place with the
# relocated LAIS record, now carrying an updated balance. Do this first,
# while row numbers still match the ORIGINAL (pre-edit) layout.
# The leading apostrophe forces the account number to stay text (so the
# leading zero in "004230529" is preserved, matching the other Conta cells).
$ws.Cells.Item(5, 1).Value = "'004230529"
$ws.Cells.Item(5, 2).Value = "LAIS"
$ws.Cells.Item(5, 3).Value = 45901.8

# Delete rows from bottom to top so earlier row numbers remain valid as we go.
# Row numbers below refer to the ORIGINAL (pre-edit) layout of the sheet.
#   146: 004230529 LAIS     184.92      (old LAIS entry, removed from here;
#                                        LAIS reappears above with a new balance)
#   16:  004974089 CELIA    5000
#   15:  005002390 LUCIANO  6000
#   14:  004211807 EDINARDO 6593.84
#   11:  002786022 PAULO    10000
#   9:   005101676 ELENI    14973.76
#   7:   004450724 ASSAKO   27030.92
#   3:   005262440 BERNARDO 113578.52
#   2:   004431546 GABRIELA 169231.01
$rowsToDelete = @(146, 16, 15, 14, 11, 9, 7, 3, 2)
foreach ($r in $rowsToDelete) {
    $ws.Rows($r).Delete()
}
